# Implemented timeseries to supim file, demand file next

$wb = $excel.ActiveWorkbook

# --- SupIm sheet: extend timeseries rows 4..14, mirroring row 3 pattern ---
$supim = $wb.Worksheets.Item("SupIm")

for ($r = 4; $r -le 14; $r++) {
    $t = $r - 2
    $supim.Cells.Item($r, 1).Value = $t
    $supim.Cells.Item($r, 2).Value = 0.481
    $supim.Cells.Item($r, 3).Value = 0.3
    $supim.Cells.Item($r, 4).Value = 0.207
}

# Make SupIm the active/selected sheet with A1:D14 selected
$supim.Activate()
$supim.Range("A1:D14").Select()

$wb.Save()
